# Applies crypto list price/volume/name/link updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "42.824.65"
    "E2" = "  +0.56%  "
    "D3" = "2.542.94"
    "E3" = "  -0.09%  "
    "E4" = "  -0.18%  "
    "D5" = "303.68"
    "E5" = "  +1.64%  "
    "D6" = "97.71"
    "E6" = "  +6.08%  "
    "E7" = "  +0.44%  "
    "E8" = "  +0.06%  "
    "D9" = "0.544"
    "E9" = "  -1.00%  "
    "D10" = "36.79"
    "E10" = "  +2.76%  "
    "D11" = "0.0828"
    "E11" = "  +3.24%  "
    "B12" = "TRON"
    "C12" = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
    "D12" = "0.114"
    "E12" = "  +0.91%  "
    "B13" = "Polkadot"
    "C13" = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
    "D13" = "7.59"
    "E13" = "  -1.42%  "
    "D14" = "2.935.41"
    "E14" = "  +0.06%  "
    "D15" = "2.527.39"
    "E15" = "  -1.47%  "
    "D16" = "15.06"
    "E16" = "  +5.77%  "
    "D17" = "0.866"
    "E17" = "  -0.99%  "
    "D18" = "42.844.40"
    "D19" = "13.33"
    "E19" = "  +3.80%  "
    "D20" = "0.0₃0988"
    "E20" = "  +0.93%  "
    "D21" = "6.57"
    "E21" = "  +0.47%  "
    "D22" = "71.79"
    "E22" = "  +0.79%  "
    "D23" = "256.05"
    "E23" = "  +0.16%  "
    "D24" = "2.94"
    "E24" = "  +1.14%  "
    "E25" = "  -1.66%  "
    "D26" = "28.10"
    "E26" = "  -3.50%  "
    "E27" = "  -0.17%  "
    "E28" = "  +9.40%  "
    "D29" = "10.16"
    "E29" = "  +1.28%  "
    "D30" = "37.86"
    "E30" = "  +2.64%  "
    "D31" = "6.07"
    "E31" = "  +2.01%  "
    "D32" = "157.92"
    "E32" = "  +3.62%  "
    "D33" = "19.52"
    "E33" = "  +14.11%  "
    "E34" = "  -1.86%  "
    "D35" = "3.31"
    "E35" = "  -2.25%  "
    "D36" = "0.0797"
    "E36" = "  +0.80%  "
    "E37" = "  -4.38%  "
    "E38" = "  +1.61%  "
    "D39" = "25.73"
    "E39" = "  +7.57%  "
    "E40" = "  +0.19%  "
    "E41" = "  +29.00%  "
    "B42" = "RenderToken"
    "C42" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D42" = "3.88"
    "E42" = "  +0.36%  "
    "B43" = "NEARProtocol"
    "C43" = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
    "D43" = "3.40"
    "E43" = "  +0.42%  "
    "D44" = "2.091.68"
    "E44" = "  +0.32%  "
    "E45" = "  -1.52%  "
    "E46" = "  +0.02%  "
    "D47" = "87.56"
    "E47" = "  +3.84%  "
    "D48" = "8.92"
    "E48" = "  -2.89%  "
    "D49" = "2.794.02"
    "E49" = "  +0.14%  "
    "D50" = "74.71"
    "E50" = "  +8.61%  "
    "E51" = "  +1.62%  "
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    # Force text interpretation so numeric-looking strings (e.g. "28.10")
    # keep their exact textual representation instead of being coerced
    # into a Double (which would drop trailing zeros / change precision).
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    $cell.Style = $origStyle
}

Write-Host "Applied $($updates.Count) cell updates"
